# Updated symbol list on Fri Jan 27 19:07:06 UTC 2023 with GitHub Actions
# Re-applies the refreshed coinranking.com snapshot (prices, 1h volume %,
# hour-of-day, and a few re-ranked coin name/link pairs) onto before.xlsx.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''308.51'
$ws.Range("E2").Value = '''1.19%'
$ws.Range("G2").Value = '''19'

# Row 3
$ws.Range("D3").Value = '''36.38'
$ws.Range("E3").Value = '''1.40%'
$ws.Range("G3").Value = '''19'

# Row 4
$ws.Range("D4").Value = '''5.064'
$ws.Range("E4").Value = '''0.88%'
$ws.Range("G4").Value = '''19'

# Row 5
$ws.Range("D5").Value = '''0.08111'
$ws.Range("E5").Value = '''0.12%'
$ws.Range("G5").Value = '''19'

# Row 6
$ws.Range("D6").Value = '''2.052'
$ws.Range("E6").Value = '''6.99%'
$ws.Range("G6").Value = '''19'

# Row 7
$ws.Range("B7").Value = '''KuCoinToken'
$ws.Range("C7").Value = '''https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").Value = '''7.858'
$ws.Range("E7").Value = '''0.24%'
$ws.Range("G7").Value = '''19'

# Row 8
$ws.Range("B8").Value = '''MXToken'
$ws.Range("C8").Value = '''https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '''0.9301'
$ws.Range("E8").Value = '''-0.28%'
$ws.Range("G8").Value = '''19'

# Row 9
$ws.Range("B9").Value = '''LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = '''https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '''0.1411'
$ws.Range("E9").Value = '''9.45%'
$ws.Range("G9").Value = '''19'

# Row 10
$ws.Range("B10").Value = '''WazirX'
$ws.Range("C10").Value = '''https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1933'
$ws.Range("E10").Value = '''0.98%'
$ws.Range("G10").Value = '''19'

# Row 11
$ws.Range("B11").Value = '''MandalaExchangeToken'
$ws.Range("C11").Value = '''https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '''0.09119'
$ws.Range("E11").Value = '''-0.95%'
$ws.Range("G11").Value = '''19'

# Row 12
$ws.Range("B12").Value = '''BitrueCoin'
$ws.Range("C12").Value = '''https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '''0.03451'
$ws.Range("E12").Value = '''-0.99%'
$ws.Range("G12").Value = '''19'

# Row 13
$ws.Range("B13").Value = '''BitMartToken'
$ws.Range("C13").Value = '''https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '''0.09914'
$ws.Range("E13").Value = '''-0.09%'
$ws.Range("G13").Value = '''19'

# Row 14
$ws.Range("B14").Value = '''BitForexToken'
$ws.Range("C14").Value = '''https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '''0.001404'
$ws.Range("E14").Value = '''-1.03%'
$ws.Range("G14").Value = '''19'

# Row 15
$ws.Range("B15").Value = '''TigerCash'
$ws.Range("C15").Value = '''https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '''0.006276'
$ws.Range("E15").Value = '''-5.57%'
$ws.Range("G15").Value = '''19'

# Row 16
$ws.Range("B16").Value = '''LEO'
$ws.Range("C16").Value = '''https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '''3.843'
$ws.Range("E16").Value = '''6.35%'
$ws.Range("G16").Value = '''19'

# Row 17
$ws.Range("B17").Value = '''GateToken'
$ws.Range("C17").Value = '''https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '''4.154'
$ws.Range("E17").Value = '''0.16%'
$ws.Range("G17").Value = '''19'

# Row 18
$ws.Range("D18").Value = '''3.462'
$ws.Range("E18").Value = '''8.08%'
$ws.Range("G18").Value = '''19'

# Row 19
$ws.Range("D19").Value = '''0.3438'
$ws.Range("E19").Value = '''-0.26%'
$ws.Range("G19").Value = '''19'

# Row 20
$ws.Range("D20").Value = '''0.1282'
$ws.Range("E20").Value = '''-4.20%'
$ws.Range("G20").Value = '''19'

# Row 21
$ws.Range("D21").Value = '''4.800'
$ws.Range("E21").Value = '''-7.83%'
$ws.Range("G21").Value = '''19'

# Row 22
$ws.Range("G22").Value = '''19'

# Row 23
$ws.Range("D23").Value = '''0.04389'
$ws.Range("E23").Value = '''-0.63%'
$ws.Range("G23").Value = '''19'

# Row 24
$ws.Range("D24").Value = '''0.001234'
$ws.Range("E24").Value = '''-0.12%'
$ws.Range("G24").Value = '''19'

# Row 25
$ws.Range("E25").Value = '''4.00%'
$ws.Range("G25").Value = '''19'

# Row 26
$ws.Range("G26").Value = '''19'

# Row 27
$ws.Range("D27").Value = '''0.0001300'
$ws.Range("E27").Value = '''-0.27%'
$ws.Range("G27").Value = '''19'

# Row 28
$ws.Range("G28").Value = '''19'

# Row 29
$ws.Range("G29").Value = '''19'

# Row 30
$ws.Range("G30").Value = '''19'

# Row 31
$ws.Range("G31").Value = '''19'

# Row 32
$ws.Range("G32").Value = '''19'

# Row 33
$ws.Range("G33").Value = '''19'

# Row 34
$ws.Range("G34").Value = '''19'

# Row 35
$ws.Range("G35").Value = '''19'

# Row 36
$ws.Range("G36").Value = '''19'

# Row 37
$ws.Range("G37").Value = '''19'

# Row 38
$ws.Range("G38").Value = '''19'

# Row 39
$ws.Range("D39").Value = '''0.02031'
$ws.Range("E39").Value = '''2.59%'
$ws.Range("G39").Value = '''19'

# Row 40
$ws.Range("D40").Value = '''0.05168'
$ws.Range("E40").Value = '''-0.07%'
$ws.Range("G40").Value = '''19'

# Row 41
$ws.Range("D41").Value = '''0.007486'
$ws.Range("E41").Value = '''-1.94%'
$ws.Range("G41").Value = '''19'

# Row 42
$ws.Range("D42").Value = '''0.01012'
$ws.Range("E42").Value = '''0.32%'
$ws.Range("G42").Value = '''19'

# Row 43
$ws.Range("D43").Value = '''0.1371'
$ws.Range("E43").Value = '''0.07%'
$ws.Range("G43").Value = '''19'

# Row 44
$ws.Range("D44").Value = '''0.002121'
$ws.Range("E44").Value = '''0.68%'
$ws.Range("G44").Value = '''19'

# Row 45
$ws.Range("D45").Value = '''0.009982'
$ws.Range("E45").Value = '''-6.83%'
$ws.Range("G45").Value = '''19'

# Row 46
$ws.Range("D46").Value = '''0.00006280'
$ws.Range("E46").Value = '''-0.64%'
$ws.Range("G46").Value = '''19'

# Row 47
$ws.Range("E47").Value = '''-0.32%'
$ws.Range("G47").Value = '''19'

# Row 48
$ws.Range("D48").Value = '''63.78'
$ws.Range("E48").Value = '''-1.81%'
$ws.Range("G48").Value = '''19'

# Row 49
$ws.Range("E49").Value = '''-22.02%'
$ws.Range("G49").Value = '''19'

# Row 50
$ws.Range("E50").Value = '''-0.32%'
$ws.Range("G50").Value = '''19'

# Row 51
$ws.Range("E51").Value = '''-0.32%'
$ws.Range("G51").Value = '''19'
